# microquiz11m-afternoon.pptx — "add afternoon header"
#
# 1) Content Placeholder 2 (the "Microquiz ..." header):
#      - reposition/resize the textbox
#      - change "Apr. 22, 2013" -> "Apr. 22, 2013 -afternoon"
# 2) TextBox 3 (equation term description):
#      - no textual change, just run-merges that happen for free when we
#        touch the (already-identical) substrings "w+x+y+z" and "a,b,c,d,e"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape 2: "Content Placeholder 2"
# ---------------------------------------------------------------------
$header = $s.Shapes.Item(2)

# Move + resize the placeholder. The target EMU values are not exact
# multiples of 1/12700" in floating point, so nudge each literal to the
# float32 points value that this host's Shape.Left/Top/Width/Height
# storage rounds back down to the exact target EMU.
$header.Left   = 111.39378356933594   # -> 1414701 EMU
$header.Top    = 14.287480354309082   # ->  181451 EMU
$header.Width  = 498.0950622558594    # -> 6325807 EMU
$header.Height = 55.74142074584961    # ->  707916 EMU

$tr = $header.TextFrame.TextRange
# "Microquiz Apr. 22, 2013"
#  123456789012345678901234
# chars 10-19 == " Apr. 22, "; chars 20-23 == "2013"
$tr.Characters(10, 10).Text = " Apr. 22, "
$tr.Characters(20, 4).Text  = "2013 -afternoon"

# ---------------------------------------------------------------------
# Shape 3: "TextBox 3"
# ---------------------------------------------------------------------
$term = $s.Shapes.Item(3)
$termTr = $term.TextFrame.TextRange

$wxyz = $termTr.Find("w+x+y+z")
$termTr.Characters($wxyz.Start, $wxyz.Length).Text = "w+x+y+z"

$abcde = $termTr.Find("a,b,c,d,e")
$termTr.Characters($abcde.Start, $abcde.Length).Text = "a,b,c,d,e"
